# Updates the cryptos list worksheet with refreshed prices / 1h volume
# percentages (and re-orders the Toncoin/Cosmos entries), matching the
# "Updated cryptos list ... with GitHub Actions" commit.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column D holds plain-looking decimal numbers (e.g. "309.41") that must stay
# TEXT (they are stored as inline/shared strings in the workbook, not as
# numbers - some prices use "." as a thousands separator e.g. "43.572.92").
# Assigning such strings straight to Range.Value lets Excel auto-convert them
# to real floating point numbers (losing formatting/precision, e.g. "0.510"
# -> 0.51) and also changes the cell's storage type. To keep them as text
# without touching the cell's style, we build the text via a formula that
# evaluates to a string, then paste-special just the resulting value back
# over itself (values only), which collapses it to a plain text cell.
function Set-TextNumber($cellref, $val) {
    $ws.Range($cellref).Formula = "=""$val"""
    $ws.Range($cellref).Copy()
    $ws.Range($cellref).PasteSpecial(-4163)
}

$ws.Range('D2').Value = '43.572.92'
$ws.Range('E2').Value = '  +0.98%  '
$ws.Range('D3').Value = '2.372.53'
$ws.Range('E3').Value = '  +2.87%  '
$ws.Range('E4').Value = '  +0.00%  '
Set-TextNumber 'D5' '309.41'
$ws.Range('E5').Value = '  -0.36%  '
Set-TextNumber 'D6' '104.71'
$ws.Range('E6').Value = '  +3.42%  '
Set-TextNumber 'D7' '0.510'
$ws.Range('E7').Value = '  -5.19%  '
$ws.Range('E8').Value = '  +0.00%  '
$ws.Range('E9').Value = '  -0.84%  '
Set-TextNumber 'D10' '35.95'
$ws.Range('E10').Value = '  -0.09%  '
Set-TextNumber 'D11' '53.41'
$ws.Range('E11').Value = '  +2.19%  '
Set-TextNumber 'D12' '0.0812'
$ws.Range('E12').Value = '  -0.86%  '
$ws.Range('E13').Value = '  -0.84%  '
Set-TextNumber 'D14' '6.99'
$ws.Range('E14').Value = '  -2.35%  '
$ws.Range('D15').Value = '2.739.28'
$ws.Range('E15').Value = '  +2.79%  '
Set-TextNumber 'D16' '15.59'
$ws.Range('E16').Value = '  +3.93%  '
$ws.Range('D17').Value = '2.376.10'
$ws.Range('E17').Value = '  +3.26%  '
$ws.Range('E18').Value = '  +0.02%  '
$ws.Range('D19').Value = '43.534.33'
$ws.Range('E19').Value = '  +1.04%  '
Set-TextNumber 'D20' '6.30'
$ws.Range('E20').Value = '  +3.39%  '
Set-TextNumber 'D21' '11.88'
$ws.Range('E21').Value = '  -5.16%  '
$ws.Range('E22').Value = '  -0.69%  '
Set-TextNumber 'D23' '68.37'
$ws.Range('E23').Value = '  -0.37%  '
Set-TextNumber 'D24' '240.82'
$ws.Range('E24').Value = '  -0.07%  '
$ws.Range('E25').Value = '  +1.92%  '
$ws.Range('E26').Value = '  -0.54%  '
$ws.Range('E27').Value = '  +0.16%  '
Set-TextNumber 'D28' '25.79'
$ws.Range('E28').Value = '  +4.01%  '
$ws.Range('E29').Value = '  -3.46%  '
Set-TextNumber 'D30' '36.47'
$ws.Range('E30').Value = '  -3.13%  '
$ws.Range('B31').Value = 'Cosmos'
$ws.Range('C31').Value = 'https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom'
Set-TextNumber 'D31' '9.53'
$ws.Range('E31').Value = '  -1.35%  '
$ws.Range('B32').Value = 'Toncoin'
$ws.Range('C32').Value = 'https://coinranking.com/coin/67YlI0K1b+toncoin-ton'
Set-TextNumber 'D32' '2.11'
$ws.Range('E32').Value = '  -0.09%  '
Set-TextNumber 'D33' '160.77'
$ws.Range('E33').Value = '  -4.12%  '
$ws.Range('E34').Value = '  -1.30%  '
$ws.Range('E35').Value = '  -0.01%  '
$ws.Range('E36').Value = '  +3.11%  '
Set-TextNumber 'D37' '2.53'
$ws.Range('E37').Value = '  +5.81%  '
Set-TextNumber 'D38' '3.11'
$ws.Range('E38').Value = '  -0.42%  '
$ws.Range('E39').Value = '  -0.38%  '
Set-TextNumber 'D40' '4.64'
$ws.Range('E40').Value = '  +8.11%  '
$ws.Range('E41').Value = '  +5.37%  '
Set-TextNumber 'D42' '0.106'
$ws.Range('E42').Value = '  -2.07%  '
Set-TextNumber 'D43' '0.114'
$ws.Range('E43').Value = '  -2.05%  '
Set-TextNumber 'D44' '2.58'
$ws.Range('E44').Value = '  +12.46%  '
$ws.Range('D45').Value = '2.035.34'
$ws.Range('E45').Value = '  +2.60%  '
Set-TextNumber 'D46' '19.62'
$ws.Range('E46').Value = '  +2.28%  '
$ws.Range('E47').Value = '  +0.33%  '
$ws.Range('E48').Value = '  +3.44%  '
$ws.Range('E49').Value = '  +7.18%  '
Set-TextNumber 'D50' '58.07'
$ws.Range('E50').Value = '  +4.31%  '
$ws.Range('E51').Value = '  -0.01%  '

$excel.CutCopyMode = 0

